$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 26 de Marzo de 2020 a las 18:42'
$ws.Cells.Item(6, 2).Value = 75665
$ws.Cells.Item(6, 3).Value = 7454
$ws.Cells.Item(6, 5).Value = 72702
$ws.Cells.Item(6, 6).Value = 2096
$ws.Cells.Item(6, 7).Value = 73
$ws.Cells.Item(6, 8).Value = 1100
$ws.Cells.Item(15, 2).Value = 6703
$ws.Cells.Item(15, 3).Value = 1115
$ws.Cells.Item(15, 5).Value = 6542
$ws.Cells.Item(19, 2).Value = 3316
$ws.Cells.Item(19, 3).Value = 232
$ws.Cells.Item(19, 5).Value = 3296
$ws.Cells.Item(26, 1).Value = 'Chequia'
$ws.Cells.Item(26, 2).Value = 1925
$ws.Cells.Item(26, 3).Value = 271
$ws.Cells.Item(26, 4).Value = 10
$ws.Cells.Item(26, 5).Value = 1906
$ws.Cells.Item(26, 6).Value = 34
$ws.Cells.Item(26, 7).Value = 3
$ws.Cells.Item(26, 8).Value = 9
$ws.Cells.Item(27, 1).Value = 'Dinamarca'
$ws.Cells.Item(27, 2).Value = 1877
$ws.Cells.Item(27, 3).Value = 153
$ws.Cells.Item(27, 4).Value = 1
$ws.Cells.Item(27, 5).Value = 1835
$ws.Cells.Item(27, 6).Value = 94
$ws.Cells.Item(27, 7).Value = 7
$ws.Cells.Item(27, 8).Value = 41
$ws.Cells.Item(48, 1).Value = 'Peru'
$ws.Cells.Item(48, 2).Value = 580
$ws.Cells.Item(48, 3).Value = 100
$ws.Cells.Item(48, 4).Value = 14
$ws.Cells.Item(48, 5).Value = 557
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 9
$ws.Cells.Item(49, 1).Value = 'Eslovenia'
$ws.Cells.Item(49, 2).Value = 562
$ws.Cells.Item(49, 3).Value = 34
$ws.Cells.Item(49, 4).Value = 10
$ws.Cells.Item(49, 5).Value = 546
$ws.Cells.Item(49, 6).Value = 14
$ws.Cells.Item(49, 7).Value = 1
$ws.Cells.Item(49, 8).Value = 6
$ws.Cells.Item(50, 1).Value = 'Panama'
$ws.Cells.Item(50, 2).Value = 558
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 2
$ws.Cells.Item(50, 5).Value = 548
$ws.Cells.Item(50, 6).Value = 20
$ws.Cells.Item(50, 8).Value = 8
$ws.Cells.Item(51, 1).Value = 'Estonia'
$ws.Cells.Item(51, 2).Value = 538
$ws.Cells.Item(51, 3).Value = 134
$ws.Cells.Item(51, 4).Value = 8
$ws.Cells.Item(51, 5).Value = 529
$ws.Cells.Item(51, 8).Value = 1
$ws.Cells.Item(52, 1).Value = 'Catar'
$ws.Cells.Item(52, 2).Value = 537
$ws.Cells.Item(52, 4).Value = 41
$ws.Cells.Item(52, 5).Value = 496
$ws.Cells.Item(52, 6).Value = 6
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(53, 1).Value = 'Argentina'
$ws.Cells.Item(53, 2).Value = 502
$ws.Cells.Item(53, 3).Value = 0
$ws.Cells.Item(53, 4).Value = 52
$ws.Cells.Item(53, 5).Value = 442
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 8
$ws.Cells.Item(54, 1).Value = 'Croacia'
$ws.Cells.Item(54, 2).Value = 495
$ws.Cells.Item(54, 3).Value = 53
$ws.Cells.Item(54, 4).Value = 22
$ws.Cells.Item(54, 5).Value = 471
$ws.Cells.Item(54, 6).Value = 14
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 2
$ws.Cells.Item(55, 1).Value = 'Republica Dominicana'
$ws.Cells.Item(55, 2).Value = 488
$ws.Cells.Item(55, 3).Value = 96
$ws.Cells.Item(55, 4).Value = 3
$ws.Cells.Item(55, 5).Value = 475
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 8).Value = 10
$ws.Cells.Item(58, 1).Value = 'Barein'
$ws.Cells.Item(58, 2).Value = 458
$ws.Cells.Item(58, 3).Value = 39
$ws.Cells.Item(58, 4).Value = 204
$ws.Cells.Item(58, 5).Value = 250
$ws.Cells.Item(58, 6).Value = 1
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 4
$ws.Cells.Item(59, 1).Value = 'Serbia'
$ws.Cells.Item(59, 3).Value = 73
$ws.Cells.Item(59, 4).Value = 15
$ws.Cells.Item(59, 5).Value = 435
$ws.Cells.Item(59, 6).Value = 21
$ws.Cells.Item(59, 7).Value = 3
$ws.Cells.Item(59, 8).Value = 7
$ws.Cells.Item(64, 4).Value = 29
$ws.Cells.Item(64, 5).Value = 313
$ws.Cells.Item(69, 1).Value = 'Marruecos'
$ws.Cells.Item(69, 2).Value = 275
$ws.Cells.Item(69, 3).Value = 50
$ws.Cells.Item(69, 5).Value = 257
$ws.Cells.Item(69, 6).Value = 1
$ws.Cells.Item(69, 7).Value = 4
$ws.Cells.Item(69, 8).Value = 10
$ws.Cells.Item(70, 1).Value = 'Bulgaria'
$ws.Cells.Item(70, 2).Value = 264
$ws.Cells.Item(70, 3).Value = 22
$ws.Cells.Item(70, 4).Value = 8
$ws.Cells.Item(70, 5).Value = 253
$ws.Cells.Item(70, 6).Value = 8
$ws.Cells.Item(70, 8).Value = 3
$ws.Cells.Item(71, 1).Value = 'Hungria'
$ws.Cells.Item(71, 2).Value = 261
$ws.Cells.Item(71, 3).Value = 35
$ws.Cells.Item(71, 4).Value = 28
$ws.Cells.Item(71, 5).Value = 223
$ws.Cells.Item(71, 6).Value = 6
$ws.Cells.Item(71, 8).Value = 10
$ws.Cells.Item(72, 1).Value = 'Taiwan'
$ws.Cells.Item(72, 2).Value = 252
$ws.Cells.Item(72, 3).Value = 17
$ws.Cells.Item(72, 4).Value = 29
$ws.Cells.Item(72, 5).Value = 221
$ws.Cells.Item(72, 8).Value = 2
$ws.Cells.Item(73, 1).Value = 'Letonia'
$ws.Cells.Item(73, 2).Value = 244
$ws.Cells.Item(73, 3).Value = 23
$ws.Cells.Item(73, 4).Value = 1
$ws.Cells.Item(73, 5).Value = 243
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(74, 1).Value = 'Eslovaquia'
$ws.Cells.Item(74, 2).Value = 226
$ws.Cells.Item(74, 3).Value = 10
$ws.Cells.Item(74, 4).Value = 2
$ws.Cells.Item(74, 5).Value = 224
$ws.Cells.Item(74, 6).Value = 2
$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(89, 1).Value = 'Republica de Chipre'
$ws.Cells.Item(89, 2).Value = 146
$ws.Cells.Item(89, 3).Value = 14
$ws.Cells.Item(89, 4).Value = 4
$ws.Cells.Item(89, 5).Value = 139
$ws.Cells.Item(89, 6).Value = 3
$ws.Cells.Item(89, 8).Value = 3
$ws.Cells.Item(90, 1).Value = 'Islas Feroe'
$ws.Cells.Item(90, 2).Value = 140
$ws.Cells.Item(90, 3).Value = 8
$ws.Cells.Item(90, 4).Value = 47
$ws.Cells.Item(90, 5).Value = 93
$ws.Cells.Item(91, 1).Value = 'Reunion'
$ws.Cells.Item(91, 2).Value = 135
$ws.Cells.Item(91, 3).Value = 24
$ws.Cells.Item(91, 4).Value = 1
$ws.Cells.Item(91, 5).Value = 134
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(92, 1).Value = 'Malta'
$ws.Cells.Item(92, 2).Value = 134
$ws.Cells.Item(92, 3).Value = 5
$ws.Cells.Item(92, 4).Value = 2
$ws.Cells.Item(92, 5).Value = 132
$ws.Cells.Item(92, 6).Value = 1
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(93, 1).Value = 'Ghana'
$ws.Cells.Item(93, 3).Value = 64
$ws.Cells.Item(93, 4).Value = 1
$ws.Cells.Item(93, 5).Value = 127
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 8).Value = 4

$wb.Save()
